# Insert a new worksheet "Most Rel Authors" right before "Most Cited Docs",
# mirroring the author-metrics export added by run_bibliometric_analysis.R.

$wb = $excel.ActiveWorkbook

$mostCitedDocs = $wb.Worksheets.Item("Most Cited Docs")

$newSheet = $wb.Worksheets.Add($mostCitedDocs)
$newSheet.Name = "Most Rel Authors"

$headers = @("Authors       ", "Articles", "Authors       ", "Articles Fractionalized")

$data = @(
  @("CHANG J-S          ", "5", "ZHU L         ", "2.167"),
  @("CHEN W-H           ", "4", "MALCATA FX    ", "1.792"),
  @("MALCATA FX         ", "4", "KUMAR A       ", "1.125"),
  @("ZHU L              ", "4", "FALFUSHYNSKA H", "1.000"),
  @("HO S-H             ", "3", "NARAYANAN M   ", "1.000"),
  @("LAM MK             ", "3", "OLTRA C       ", "1.000"),
  @("LIM JUN WEI        ", "3", "SINGH J       ", "1.000"),
  @("LIU J              ", "3", "VAN BEILEN JB ", "1.000"),
  @("SHOW PAU LOKE      ", "3", "SEN R         ", "0.833"),
  @("SHOW PL            ", "3", "CHANG J-S     ", "0.787"),
  @("WANG J             ", "3", "LAM MK        ", "0.754"),
  @("WANG Y             ", "3", "SINGH B       ", "0.750"),
  @("ABOMOHRA ABDELFATAH", "2", "AMARO HM      ", "0.667"),
  @("ABREU M            ", "2", "CHEN Y        ", "0.667"),
  @("AMARO HM           ", "2", "MISHRA S      ", "0.667"),
  @("CAETANO N          ", "2", "SIALVE B      ", "0.667"),
  @("CALIJURI ML        ", "2", "VERMA AK      ", "0.667"),
  @("CHEN H             ", "2", "KUMAR L       ", "0.625"),
  @("CHEN Y             ", "2", "HO S-H        ", "0.587"),
  @("CHEW KW            ", "2", "CAETANO N     ", "0.583")
)

# All data on this sheet is exported as plain text (matches the rest of the
# workbook, where numeric-looking values are also exported as text), so
# force the Text number format on the whole used range before writing any
# values - this keeps "5", "2.167", etc. as strings instead of numbers.
$usedRange = $newSheet.Range("A1:D21")
$usedRange.NumberFormat = "@"

for ($col = 1; $col -le 4; $col++) {
  $newSheet.Cells.Item(1, $col).Value = $headers[$col - 1]
}

$row = 2
foreach ($r in $data) {
  $newSheet.Cells.Item($row, 1).Value = $r[0]
  $newSheet.Cells.Item($row, 2).Value = $r[1]
  $newSheet.Cells.Item($row, 3).Value = $r[2]
  $newSheet.Cells.Item($row, 4).Value = $r[3]
  $row++
}

# Drop the temporary Text number format again so the cells end up with no
# explicit style (matching the other, unstyled, data sheets) while keeping
# the values stored as text.
$usedRange.ClearFormats()

# Header row (bold, centered - same style as the other sheets' headers)
$headerRange = $newSheet.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
